# Auto-generated Excel COM-interop script to apply diff changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect before editing, then reapply protection afterwards
$ws.Unprotect()

# Update the confidential disclosure text date (2021-04-30 -> 2021-05-03)
$ws.Range("A80").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-03 for illustrative purposes only and are subject to change."

# Update weight (D) and return (E) values for each holding row
$ws.Range("D2").Value = 0.065794026464722
$ws.Range("E2").Value = 0.00821542674577791
$ws.Range("D3").Value = 0.04131901850360174
$ws.Range("E3").Value = -0.02334012032000743
$ws.Range("D4").Value = 0.03305573912263907
$ws.Range("E4").Value = -0.001268934887778461
$ws.Range("D5").Value = 0.03061503065066056
$ws.Range("E5").Value = -0.004067477307758249
$ws.Range("D6").Value = 0.02804514885656387
$ws.Range("E6").Value = -0.004427448480985752
$ws.Range("D7").Value = 0.02443806725092041
$ws.Range("E7").Value = -0.002925687536570964
$ws.Range("D8").Value = 0.1796112939587997
$ws.Range("E8").Value = 0.03001017293997976
$ws.Range("D9").Value = 0.02391617048606467
$ws.Range("E9").Value = 0.01523996804522842
$ws.Range("D10").Value = 0.02225832702196665
$ws.Range("E10").Value = 0.01184230250337293
$ws.Range("D11").Value = 0.02246726050116151
$ws.Range("E11").Value = 0.001202212070209097
$ws.Range("D12").Value = 0.02068902209876699
$ws.Range("E12").Value = -0.002741640683797564
$ws.Range("D13").Value = 0.01964077979535402
$ws.Range("E13").Value = 0.0007401924500372026
$ws.Range("D14").Value = 0.01678430958825527
$ws.Range("E14").Value = 0.005107051659791972
$ws.Range("D15").Value = 0.01723534374343731
$ws.Range("E15").Value = -0.01534887579855626
$ws.Range("D16").Value = 0.01544360013516328
$ws.Range("E16").Value = -0.01160493827160491
$ws.Range("D17").Value = 0.01444921977049332
$ws.Range("E17").Value = 0.01579586877278261
$ws.Range("D18").Value = 0.01431551823228612
$ws.Range("E18").Value = 0.01130688124306323
$ws.Range("D19").Value = 0.01420382223542378
$ws.Range("E19").Value = -0.007690414667158807
$ws.Range("D20").Value = 0.01273239033365683
$ws.Range("E20").Value = 0.02760307477288615
$ws.Range("D21").Value = 0.0123516626914358
$ws.Range("E21").Value = 0.002865329512893977
$ws.Range("D22").Value = 0.0128438480033566
$ws.Range("E22").Value = 0.02852017937219742
$ws.Range("D23").Value = 0.01214070343135902
$ws.Range("E23").Value = -0.006962257237083214
$ws.Range("D24").Value = 0.01267264965823686
$ws.Range("E24").Value = 0.0178034102306921
$ws.Range("D25").Value = 0.01157090290682476
$ws.Range("E25").Value = 0.0203911391231808
$ws.Range("D26").Value = 0.009572648812175436
$ws.Range("E26").Value = -0.01219937260369475
$ws.Range("D27").Value = 0.009858005868170839
$ws.Range("E27").Value = -0.005173664275928802
$ws.Range("D28").Value = 0.01018522906239274
$ws.Range("E28").Value = 0.01334539696901138
$ws.Range("D29").Value = 0.01044222518604116
$ws.Range("E29").Value = 0.0001293326435591879
$ws.Range("D30").Value = 0.01017498099440181
$ws.Range("E30").Value = -0.000468457214241158
$ws.Range("D31").Value = 0.008912125081085851
$ws.Range("E31").Value = -0.004693203545976132
$ws.Range("D32").Value = 0.0102595871371176
$ws.Range("E32").Value = 0.007479964381121862
$ws.Range("D33").Value = 0.009411499929077808
$ws.Range("E33").Value = 0.002249524139124492
$ws.Range("D34").Value = 0.009005430165235796
$ws.Range("E34").Value = 0.009262689885142672
$ws.Range("D35").Value = 0.009154146314685516
$ws.Range("E35").Value = 0.01345135815325871
$ws.Range("D36").Value = 0.008585497704770399
$ws.Range("E36").Value = -0.004580279813457766
$ws.Range("D37").Value = 0.008516025336878424
$ws.Range("E37").Value = 0.01000489750227374
$ws.Range("D38").Value = 0.008453941111026417
$ws.Range("E38").Value = -0.03459066305818692
$ws.Range("D39").Value = 0.008597255178124331
$ws.Range("E39").Value = 0.03053040103492877
$ws.Range("D40").Value = 0.007778720538244419
$ws.Range("E40").Value = 0.02386727466770155
$ws.Range("D41").Value = 0.00724832343792334
$ws.Range("E41").Value = -0.01328364752301625
$ws.Range("D42").Value = 0.007718622372080577
$ws.Range("E42").Value = -0.00987031700288199
$ws.Range("D43").Value = 0.007827220115829385
$ws.Range("E43").Value = 0.02031930333817122
$ws.Range("D44").Value = 0.007374954603640678
$ws.Range("E44").Value = 0.003091539737596127
$ws.Range("D45").Value = 0.007471238777322868
$ws.Range("E45").Value = -0.006060863832592633
$ws.Range("D46").Value = 0.007729426536784188
$ws.Range("E46").Value = -0.001726689689195871
$ws.Range("D47").Value = 0.00710214944487447
$ws.Range("E47").Value = 0.01986577181208049
$ws.Range("D48").Value = 0.00717348870887332
$ws.Range("E48").Value = 0.006622516556291647
$ws.Range("D49").Value = 0.006821519210936896
$ws.Range("E49").Value = 0.0107432963577605
$ws.Range("D50").Value = 0.006676338247732109
$ws.Range("E50").Value = 0.01688481675392639
$ws.Range("D51").Value = 0.006342561056539265
$ws.Range("E51").Value = 0.0002630309938187469
$ws.Range("D52").Value = 0.006508714810051066
$ws.Range("E52").Value = 0.007713902111558646
$ws.Range("D53").Value = 0.005591155234119298
$ws.Range("E53").Value = -0.007459505541347111
$ws.Range("D54").Value = 0.005807874067291754
$ws.Range("E54").Value = 0.01723477594791256
$ws.Range("D55").Value = 0.005856294202489193
$ws.Range("E55").Value = 0.02604537592837519
$ws.Range("D56").Value = 0.005660865929173854
$ws.Range("E56").Value = 0.007037153983791233
$ws.Range("D57").Value = 0.006478050048465815
$ws.Range("E57").Value = 0.04091042872559592
$ws.Range("D58").Value = 0.005553976196756867
$ws.Range("E58").Value = 0.000343288705801692
$ws.Range("D59").Value = 0.005093528118652925
$ws.Range("E59").Value = 0.006597417181358756
$ws.Range("D60").Value = 0.004919708174744812
$ws.Range("E60").Value = 0.004392197390518104
$ws.Range("D61").Value = 0.004532744305102937
$ws.Range("E61").Value = 0.03981106612685581
$ws.Range("D62").Value = 0.004583428548344883
$ws.Range("E62").Value = 0.02565213623364238
$ws.Range("D63").Value = 0.004236503642016395
$ws.Range("E63").Value = 0.01065106510651082
$ws.Range("D64").Value = 0.003967035063526304
$ws.Range("E64").Value = 0.01698173662287727
$ws.Range("D65").Value = 0.003953212088096682
$ws.Range("E65").Value = 0.01732245488525375
$ws.Range("D66").Value = 0.003691687748359246
$ws.Range("E66").Value = 0.01355713363460298
$ws.Range("D67").Value = 0.003833015755768997
$ws.Range("E67").Value = -0.001430081452465481
$ws.Range("D68").Value = 0.003459437928425348
$ws.Range("E68").Value = 0.04905101443284776
$ws.Range("D69").Value = 0.00353709286223256
$ws.Range("E69").Value = 0.01203845117240143
$ws.Range("D70").Value = 0.003173882266461128
$ws.Range("E70").Value = -0.004955947136563776
$ws.Range("D71").Value = 0.003125541573651216
$ws.Range("E71").Value = -0.01220023637958001
$ws.Range("D72").Value = 0.002384781031159807
$ws.Range("E72").Value = -0.0115093773943169
$ws.Range("D73").Value = 0.00203956413616314
$ws.Range("E73").Value = -0.008491245837147332
$ws.Range("D74").Value = 0.002019187163762577
$ws.Range("E74").Value = -0.00959987410001184
$ws.Range("D75").Value = 0.001554846408669834
$ws.Range("E75").Value = -0.001072961373390635
$ws.Range("D76").Value = 0.001450856323397566
$ws.Range("E76").Value = 0.04845863220719493
$ws.Range("E77").Value = 0.008335607702685932

# Restore worksheet protection to its original (protected) state
$ws.Protect()
